# Applies the F-column ("想去人数" / interest count) updates captured in the
# commit diff, across all four worksheets of the workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 37837
$ws.Range("F4").Value = 641
$ws.Range("F5").Value = 789
$ws.Range("F6").Value = 487
$ws.Range("F10").Value = 100
$ws.Range("F11").Value = 736
$ws.Range("F12").Value = 569
$ws.Range("F13").Value = 71
$ws.Range("F16").Value = 672
$ws.Range("F17").Value = 187
$ws.Range("F18").Value = 478
$ws.Range("F20").Value = 1181
$ws.Range("F21").Value = 96
$ws.Range("F22").Value = 857
$ws.Range("F23").Value = 2572
$ws.Range("F24").Value = 1054
$ws.Range("F25").Value = 573
$ws.Range("F26").Value = 111
$ws.Range("F27").Value = 1171
$ws.Range("F29").Value = 812
$ws.Range("F31").Value = 1172

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 425
$ws.Range("F4").Value = 335
$ws.Range("F9").Value = 144
$ws.Range("F12").Value = 12

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 652

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 652
$ws.Range("F3").Value = 37837
$ws.Range("F5").Value = 641
$ws.Range("F6").Value = 789
$ws.Range("F7").Value = 487
$ws.Range("F11").Value = 425
$ws.Range("F12").Value = 335
$ws.Range("F16").Value = 100
$ws.Range("F17").Value = 736
$ws.Range("F18").Value = 569
$ws.Range("F19").Value = 71
$ws.Range("F23").Value = 144
$ws.Range("F27").Value = 672
$ws.Range("F28").Value = 187
$ws.Range("F29").Value = 478
$ws.Range("F31").Value = 1181
$ws.Range("F32").Value = 96
$ws.Range("F33").Value = 858
$ws.Range("F34").Value = 2572
$ws.Range("F35").Value = 1054
$ws.Range("F36").Value = 573
$ws.Range("F37").Value = 111
$ws.Range("F38").Value = 1171
$ws.Range("F40").Value = 12
$ws.Range("F41").Value = 812
$ws.Range("F43").Value = 1172
